$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 18

# Copy the formatting (style) of the previous log row down onto the new row
# so the appended entries look the same as the rest of the log (centered
# alignment, etc. -> style index 3 in the original sheet).
$srcRange = $ws.Range("A17:H17")
$dstRange = $ws.Range("A18:H18")
$srcRange.Copy()
$dstRange.PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new run-log entry.
$ws.Cells.Item($row, 1).Value = "2025-08-16 03:56:35 UTC"
$ws.Cells.Item($row, 2).Value = "2025-08-16 09:26:35 IST"
$ws.Cells.Item($row, 3).Value = "SKIPPED"
$ws.Cells.Item($row, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($row, 5).Value = "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf"
$ws.Cells.Item($row, 6).Value = ""
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = ""
